# Coeurl_Profits market-data refresh (scheduled runner)
# Source data is plain cached values (currentAveragePrice*, LevePrice*, LeveProfit*)
# scraped from the Universalis API per-leve-item row -- no formulas involved, so this
# script just re-pokes the literal H:N values that the scraper refreshed, row by row,
# per job/sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 14285954
$ws.Range("I33").Value = 16666930
$ws.Range("K33").Value = 16666930
$ws.Range("M33").Value = -16666701
$ws.Range("H40").Value = 1829.7
$ws.Range("J40").Value = 1829.7
$ws.Range("L40").Value = 1829.7
$ws.Range("N40").Value = -2179.7
$ws.Range("H54").Value = 17796.75
$ws.Range("I54").Value = 17796.75
$ws.Range("K54").Value = 17796.75
$ws.Range("M54").Value = -17310.75
$ws.Range("H76").Value = 4482.3335
$ws.Range("I76").Value = 4482.3335
$ws.Range("K76").Value = 4482.3335
$ws.Range("M76").Value = -4167.3335
$ws.Range("H79").Value = 4482.3335
$ws.Range("I79").Value = 4482.3335
$ws.Range("K79").Value = 4482.3335
$ws.Range("M79").Value = -3390.3335
$ws.Range("H103").Value = 645
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 645
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 1935
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -3107
$ws.Range("H116").Value = 12405.286
$ws.Range("J116").Value = 9462.375
$ws.Range("L116").Value = 9462.375
$ws.Range("N116").Value = -16346.375
$ws.Range("H138").Value = 3199.9768
$ws.Range("I138").Value = 1984.8
$ws.Range("J138").Value = 3850.9644
$ws.Range("K138").Value = 5954.4
$ws.Range("L138").Value = 11552.8932
$ws.Range("M138").Value = -814.3999999999996
$ws.Range("N138").Value = -21832.8932
$ws.Range("H141").Value = 1742.125
$ws.Range("I141").Value = 1582.3182
$ws.Range("J141").Value = 3500
$ws.Range("K141").Value = 4746.9546
$ws.Range("L141").Value = 10500
$ws.Range("M141").Value = 433.0454
$ws.Range("N141").Value = -20860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 33336638
$ws.Range("I2").Value = 40003296
$ws.Range("J2").Value = 3358.2
$ws.Range("K2").Value = 40003296
$ws.Range("L2").Value = 3358.2
$ws.Range("M2").Value = -40003183
$ws.Range("N2").Value = -3584.2
$ws.Range("H45").Value = 9742.25
$ws.Range("I45").Value = 9927.134
$ws.Range("K45").Value = 9927.134
$ws.Range("M45").Value = -9550.134
$ws.Range("H106").Value = 156666.67
$ws.Range("J106").Value = 156666.67
$ws.Range("L106").Value = 156666.67
$ws.Range("N106").Value = -159190.67
$ws.Range("H116").Value = 33336638
$ws.Range("I116").Value = 40003296
$ws.Range("J116").Value = 3358.2
$ws.Range("K116").Value = 40003296
$ws.Range("L116").Value = 3358.2
$ws.Range("M116").Value = -40001002
$ws.Range("N116").Value = -7946.2
$ws.Range("H132").Value = 2875.8718
$ws.Range("J132").Value = 3602.2
$ws.Range("L132").Value = 10806.6
$ws.Range("N132").Value = -15866.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 33336638
$ws.Range("I3").Value = 40003296
$ws.Range("J3").Value = 3358.2
$ws.Range("K3").Value = 40003296
$ws.Range("L3").Value = 3358.2
$ws.Range("M3").Value = -40003182
$ws.Range("N3").Value = -3586.2
$ws.Range("H86").Value = 1959.6875
$ws.Range("I86").Value = 1959.6875
$ws.Range("K86").Value = 1959.6875
$ws.Range("M86").Value = -836.6875
$ws.Range("H89").Value = 1959.6875
$ws.Range("I89").Value = 1959.6875
$ws.Range("K89").Value = 9798.4375
$ws.Range("M89").Value = -4182.4375
$ws.Range("H96").Value = 13220.6
$ws.Range("I96").Value = 13220.6
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 13220.6
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -10474.6
$ws.Range("N96").ClearContents()
$ws.Range("H99").Value = 4153.35
$ws.Range("I99").Value = 1324.1428
$ws.Range("K99").Value = 1324.1428
$ws.Range("M99").Value = 173.8571999999999
$ws.Range("H105").Value = 2988.5386
$ws.Range("I105").Value = 1812.8334
$ws.Range("J105").Value = 5633.875
$ws.Range("K105").Value = 1812.8334
$ws.Range("L105").Value = 5633.875
$ws.Range("M105").Value = -65.83339999999998
$ws.Range("N105").Value = -9127.875
$ws.Range("H134").Value = 2028.0625
$ws.Range("I134").Value = 1646.5682
$ws.Range("J134").Value = 6224.5
$ws.Range("K134").Value = 4939.7046
$ws.Range("L134").Value = 18673.5
$ws.Range("M134").Value = -2404.7046
$ws.Range("N134").Value = -23743.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2437.5952
$ws.Range("J31").Value = 4300.8335
$ws.Range("L31").Value = 4300.8335
$ws.Range("N31").Value = -4890.8335
$ws.Range("H34").Value = 2437.5952
$ws.Range("J34").Value = 4300.8335
$ws.Range("L34").Value = 4300.8335
$ws.Range("N34").Value = -4704.8335
$ws.Range("H105").Value = 1055
$ws.Range("I105").Value = 1055
$ws.Range("K105").Value = 1055
$ws.Range("M105").Value = 692
$ws.Range("H132").Value = 3792.975
$ws.Range("I132").Value = 3558.5151
$ws.Range("K132").Value = 10675.5453
$ws.Range("M132").Value = -8145.5453
$ws.Range("H141").Value = 70241.664
$ws.Range("J141").Value = 70241.664
$ws.Range("L141").Value = 70241.664
$ws.Range("N141").Value = -80601.664

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 8999.875
$ws.Range("J116").Value = 9142.857
$ws.Range("L116").Value = 27428.571
$ws.Range("N116").Value = -34312.571
$ws.Range("H129").Value = 1069.1818
$ws.Range("J129").Value = 1553.2
$ws.Range("L129").Value = 4659.6
$ws.Range("N129").Value = -14659.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17885.111
$ws.Range("I70").Value = 14500.5
$ws.Range("J70").Value = 24654.334
$ws.Range("K70").Value = 14500.5
$ws.Range("L70").Value = 24654.334
$ws.Range("M70").Value = -14230.5
$ws.Range("N70").Value = -25194.334
$ws.Range("H73").Value = 17885.111
$ws.Range("I73").Value = 14500.5
$ws.Range("J73").Value = 24654.334
$ws.Range("K73").Value = 14500.5
$ws.Range("L73").Value = 24654.334
$ws.Range("M73").Value = -13564.5
$ws.Range("N73").Value = -26526.334
$ws.Range("H102").Value = 35715704
$ws.Range("I102").Value = 902.7143
$ws.Range("K102").Value = 902.7143
$ws.Range("M102").Value = 719.2857
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H22").Value = 1585.8695
$ws.Range("I22").Value = 1060
$ws.Range("J22").Value = 1664.75
$ws.Range("K22").Value = 1060
$ws.Range("L22").Value = 1664.75
$ws.Range("M22").Value = -765
$ws.Range("N22").Value = -2254.75
$ws.Range("H27").Value = 1585.8695
$ws.Range("I27").Value = 1060
$ws.Range("J27").Value = 1664.75
$ws.Range("K27").Value = 1060
$ws.Range("L27").Value = 1664.75
$ws.Range("M27").Value = -953
$ws.Range("N27").Value = -1878.75
$ws.Range("H42").Value = 16619.46
$ws.Range("J42").Value = 18514
$ws.Range("L42").Value = 18514
$ws.Range("N42").Value = -19640
$ws.Range("H46").Value = 3038.6
$ws.Range("J46").Value = 3273.25
$ws.Range("L46").Value = 3273.25
$ws.Range("N46").Value = -3649.25
$ws.Range("H49").Value = 16619.46
$ws.Range("J49").Value = 18514
$ws.Range("L49").Value = 18514
$ws.Range("N49").Value = -18808
$ws.Range("H55").Value = 676.1667
$ws.Range("I55").Value = 317
$ws.Range("J55").Value = 932.7143
$ws.Range("K55").Value = 317
$ws.Range("L55").Value = 932.7143
$ws.Range("M55").Value = -144
$ws.Range("N55").Value = -1278.7143

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 84365.8
$ws.Range("J46").Value = 84365.8
$ws.Range("L46").Value = 84365.8
$ws.Range("N46").Value = -84827.8
$ws.Range("H100").Value = 1380.3334
$ws.Range("I100").Value = 707.5714
$ws.Range("J100").Value = 1969
$ws.Range("K100").Value = 1415.1428
$ws.Range("L100").Value = 3938
$ws.Range("M100").Value = -874.1428000000001
$ws.Range("N100").Value = -5020
$ws.Range("H107").Value = 721.5454999999999
$ws.Range("I107").Value = 811.125
$ws.Range("J107").Value = 482.66666
$ws.Range("K107").Value = 2433.375
$ws.Range("L107").Value = 1447.99998
$ws.Range("M107").Value = -513.375
$ws.Range("N107").Value = -5287.999980000001
$ws.Range("H123").Value = 42000
$ws.Range("J123").Value = 42000
$ws.Range("L123").Value = 42000
$ws.Range("N123").Value = -51800
$ws.Range("H134").Value = 84365.8
$ws.Range("J134").Value = 84365.8
$ws.Range("L134").Value = 253097.4
$ws.Range("N134").Value = -258167.4

